$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 310, shifting existing rows 310:424 down to 311:425.
$ws.Rows("310:310").Insert()

# Populate the newly inserted row 310 with the new record's data.
$ws.Range("A310").Value = 4
$ws.Range("B310").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C310").Value = "Los Lagos"
$ws.Range("D310").Value = 44988
$ws.Range("E310").Value = 10
$ws.Range("F310").Value = 100112003
$ws.Range("G310").Value = "Ajo"
$ws.Range("H310").Value = "Chino"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 240
$ws.Range("K310").Value = 19000
$ws.Range("L310").Value = 20000
$ws.Range("M310").Value = 19500
$ws.Range("N310").Value = '$/caja 10 kilos'
$ws.Range("O310").Value = "China"
$ws.Range("P310").Value = 1950
$ws.Range("Q310").Value = 10
$ws.Range("R310").Value = "Hortaliza"
